$wb = $excel.ActiveWorkbook

# --- Overview sheet: update handback status text for both locales ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: handback completed, clear error, update handback datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-23 16:53:36"
# Error Detail cleared (leading apostrophe forces a real, typed empty string
# cell instead of simply deleting the cell); Style reset keeps no quote-prefix
# formatting behind.
$wsZhCn.Range("P2").Value = "'"
$wsZhCn.Range("P2").Style = "Normal"

# --- de-de sheet: handback completed, clear error, update handback datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-23 16:53:43"
$wsDeDe.Range("P2").Value = "'"
$wsDeDe.Range("P2").Style = "Normal"

# --- Column width adjustments ---
# Overview: columns E (zh-cn) and F (de-de) widened to fit the new status text
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

# zh-cn: Status column (C) widened; Error Detail column (P) narrowed (now empty)
$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527
$wsZhCn.Range("P1").ColumnWidth = 13.7470528738839

# de-de: Status column (C) widened; Error Detail column (P) narrowed (now empty)
$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDeDe.Range("P1").ColumnWidth = 13.7470528738839
